$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 968.26086
$ws.Range("I4").Value = 1007.7273
$ws.Range("K4").Value = 1007.7273
$ws.Range("M4").Value = -893.7273
$ws.Range("H64").Value = 46916.523
$ws.Range("J64").Value = 3624
$ws.Range("L64").Value = 3624
$ws.Range("N64").Value = -4120
$ws.Range("H67").Value = 46916.523
$ws.Range("J67").Value = 3624
$ws.Range("L67").Value = 3624
$ws.Range("N67").Value = -5340
$ws.Range("H74").Value = 5200
$ws.Range("I74").Value = 5200
$ws.Range("K74").Value = 5200
$ws.Range("M74").Value = -4264
$ws.Range("H77").Value = 5200
$ws.Range("I77").Value = 5200
$ws.Range("K77").Value = 26000
$ws.Range("M77").Value = -21320
$ws.Range("H94").Value = 83334264
$ws.Range("I94").Value = 1012.0909
$ws.Range("K94").Value = 1012.0909
$ws.Range("M94").Value = -561.0909
$ws.Range("H121").Value = 1015
$ws.Range("I121").Value = 733.3333
$ws.Range("J121").Value = 1296.6666
$ws.Range("K121").Value = 2199.9999
$ws.Range("L121").Value = 3889.9998
$ws.Range("M121").Value = -452.9998999999998
$ws.Range("N121").Value = -7383.9998
$ws.Range("H137").Value = 1279.9302
$ws.Range("I137").Value = 1099.5151
$ws.Range("J137").Value = 1875.3
$ws.Range("K137").Value = 3298.5453
$ws.Range("L137").Value = 5625.9
$ws.Range("M137").Value = -748.5453000000002
$ws.Range("N137").Value = -10725.9
$ws.Range("H138").Value = 2681.9565
$ws.Range("I138").Value = 1532.04
$ws.Range("J138").Value = 3335.318
$ws.Range("K138").Value = 4596.12
$ws.Range("L138").Value = 10005.954
$ws.Range("M138").Value = 543.8800000000001
$ws.Range("N138").Value = -20285.954

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18928.209
$ws.Range("I32").Value = 3054.6206
$ws.Range("J32").Value = 172372.89
$ws.Range("K32").Value = 3054.6206
$ws.Range("L32").Value = 172372.89
$ws.Range("M32").Value = -2767.6206
$ws.Range("N32").Value = -172946.89
$ws.Range("H61").Value = 1759.6786
$ws.Range("I61").Value = 1180.2354
$ws.Range("J61").Value = 2655.182
$ws.Range("K61").Value = 1180.2354
$ws.Range("L61").Value = 2655.182
$ws.Range("M61").Value = -968.2354
$ws.Range("N61").Value = -3079.182
$ws.Range("H136").Value = 1759.6786
$ws.Range("I136").Value = 1180.2354
$ws.Range("J136").Value = 2655.182
$ws.Range("K136").Value = 3540.7062
$ws.Range("L136").Value = 7965.545999999999
$ws.Range("M136").Value = -990.7062000000001
$ws.Range("N136").Value = -13065.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18631.889
$ws.Range("I82").Value = 3839
$ws.Range("J82").Value = 22858.428
$ws.Range("K82").Value = 3839
$ws.Range("L82").Value = 22858.428
$ws.Range("M82").Value = -3456
$ws.Range("N82").Value = -23624.428
$ws.Range("H85").Value = 18631.889
$ws.Range("I85").Value = 3839
$ws.Range("J85").Value = 22858.428
$ws.Range("K85").Value = 3839
$ws.Range("L85").Value = 22858.428
$ws.Range("M85").Value = -2513
$ws.Range("N85").Value = -25510.428
$ws.Range("H86").Value = 42380.832
$ws.Range("I86").Value = 54248.61
$ws.Range("J86").Value = 3386.7144
$ws.Range("K86").Value = 54248.61
$ws.Range("L86").Value = 3386.7144
$ws.Range("M86").Value = -53125.61
$ws.Range("N86").Value = -5632.7144
$ws.Range("H89").Value = 42380.832
$ws.Range("I89").Value = 54248.61
$ws.Range("J89").Value = 3386.7144
$ws.Range("K89").Value = 271243.05
$ws.Range("L89").Value = 16933.572
$ws.Range("M89").Value = -265627.05
$ws.Range("N89").Value = -28165.572
$ws.Range("H94").Value = 733
$ws.Range("I94").Value = 618.7143
$ws.Range("J94").Value = 999.6667
$ws.Range("K94").Value = 618.7143
$ws.Range("L94").Value = 999.6667
$ws.Range("M94").Value = -167.7143
$ws.Range("N94").Value = -1901.6667
$ws.Range("H105").Value = 135256.6
$ws.Range("I105").Value = 112664.22
$ws.Range("J105").Value = 169145.17
$ws.Range("K105").Value = 112664.22
$ws.Range("L105").Value = 169145.17
$ws.Range("M105").Value = -110917.22
$ws.Range("N105").Value = -172639.17
$ws.Range("H107").Value = 62551020
$ws.Range("I107").Value = 66720950
$ws.Range("K107").Value = 66720950
$ws.Range("M107").Value = -66719030
$ws.Range("H128").Value = 2500
$ws.Range("I128").Value = 2500
$ws.Range("K128").Value = 7500
$ws.Range("M128").Value = -5010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1985.7142
$ws.Range("I50").Value = 2120
$ws.Range("J50").Value = 1650
$ws.Range("K50").Value = 6360
$ws.Range("L50").Value = 4950
$ws.Range("M50").Value = -5879
$ws.Range("N50").Value = -5912
$ws.Range("H52").Value = 1773.6154
$ws.Range("J52").Value = 1773.6154
$ws.Range("L52").Value = 5320.8462
$ws.Range("N52").Value = -5852.8462
$ws.Range("H53").Value = 1985.7142
$ws.Range("I53").Value = 2120
$ws.Range("J53").Value = 1650
$ws.Range("K53").Value = 6360
$ws.Range("L53").Value = 4950
$ws.Range("M53").Value = -5879
$ws.Range("N53").Value = -5912
$ws.Range("H62").Value = 2733.2222
$ws.Range("J62").Value = 4280
$ws.Range("L62").Value = 12840
$ws.Range("N62").Value = -14212
$ws.Range("H65").Value = 2733.2222
$ws.Range("J65").Value = 4280
$ws.Range("L65").Value = 38520
$ws.Range("N65").Value = -45384
$ws.Range("H129").Value = 15632792
$ws.Range("I129").Value = 35714700
$ws.Range("J129").Value = 13531.111
$ws.Range("K129").Value = 107144100
$ws.Range("L129").Value = 40593.333
$ws.Range("M129").Value = -107139100
$ws.Range("N129").Value = -50593.333
$ws.Range("H137").Value = 14540034
$ws.Range("I137").Value = 87715
$ws.Range("J137").Value = 30306200
$ws.Range("K137").Value = 263145
$ws.Range("L137").Value = 90918600
$ws.Range("M137").Value = -258045
$ws.Range("N137").Value = -90928800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 19999
$ws.Range("J117").Value = 19999
$ws.Range("L117").Value = 19999
$ws.Range("N117").Value = -26883

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1032.2
$ws.Range("I22").Value = 1099.6666
$ws.Range("J22").Value = 1015.3333
$ws.Range("K22").Value = 1099.6666
$ws.Range("L22").Value = 1015.3333
$ws.Range("M22").Value = -804.6666
$ws.Range("N22").Value = -1605.3333
$ws.Range("H27").Value = 1032.2
$ws.Range("I27").Value = 1099.6666
$ws.Range("J27").Value = 1015.3333
$ws.Range("K27").Value = 1099.6666
$ws.Range("L27").Value = 1015.3333
$ws.Range("M27").Value = -992.6666
$ws.Range("N27").Value = -1229.3333
